$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.244.25'
$ws.Range("E2").Value = '  +3.53%  '
$ws.Range("D3").Value = '3.120.09'
$ws.Range("E3").Value = '  +1.78%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.39'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '625.04'
$ws.Range("E6").Value = '  +1.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.381'
$ws.Range("E7").Value = '  +2.99%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.971'
$ws.Range("E8").Value = '  +20.72%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").Value = '3.116.98'
$ws.Range("E10").Value = '  +1.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.728'
$ws.Range("E11").Value = '  +23.06%  '
$ws.Range("E12").Value = '  +5.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000255'
$ws.Range("E13").Value = '  +7.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.59'
$ws.Range("E14").Value = '  +8.50%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '91.134.80'
$ws.Range("E15").Value = '  +3.70%  '
$ws.Range("B16").Value = 'Toncoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.39'
$ws.Range("E16").Value = '  +1.97%  '
$ws.Range("D17").Value = '3.694.93'
$ws.Range("E17").Value = '  +1.76%  '
$ws.Range("D18").Value = '3.098.80'
$ws.Range("E18").Value = '  +1.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.74'
$ws.Range("E19").Value = '  +14.41%  '
$ws.Range("E20").Value = '  +10.39%  '
$ws.Range("E21").Value = '  +6.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '437.71'
$ws.Range("E22").Value = '  +4.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.78'
$ws.Range("E23").Value = '  +8.19%  '
$ws.Range("E24").Value = '  +6.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.19'
$ws.Range("E25").Value = '  +13.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.24'
$ws.Range("E26").Value = '  +4.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '86.70'
$ws.Range("E27").Value = '  +6.21%  '
$ws.Range("D28").Value = '3.290.69'
$ws.Range("E28").Value = '  +1.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.168'
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.12'
$ws.Range("E31").Value = '  +13.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.996'
$ws.Range("E32").Value = '  -8.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '526.11'
$ws.Range("E33").Value = '  +3.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.76'
$ws.Range("E34").Value = '  +5.23%  '
$ws.Range("E35").Value = '  +5.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.141'
$ws.Range("E36").Value = '  +8.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '23.73'
$ws.Range("E37").Value = '  +6.80%  '
$ws.Range("E38").Value = '  +3.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.27'
$ws.Range("E39").Value = '  +3.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.28'
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0844'
$ws.Range("E42").Value = '  +22.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.152'
$ws.Range("E43").Value = '  +14.79%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.378'
$ws.Range("E45").Value = '  +5.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.90'
$ws.Range("E46").Value = '  +5.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '147.42'
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.00'
$ws.Range("E48").Value = '  +1.57%  '
$ws.Range("E49").Value = '  +10.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '166.68'
$ws.Range("E50").Value = '  +7.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000259'
$ws.Range("E51").Value = '  +21.39%  '
